$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update product title (column A) and detail (column B) text for rows 2-18 ---
# Rows 2-15: existing product rows get replaced with new product data (content + row height).
# Rows 16-18: previously-empty placeholder rows become new data rows (content + style + height).

# Row 2
$aText = @'
【Y96031】實拍S~XL韓版高腰顯瘦白色直筒九分牛仔褲240925
'@
$bText = @'
S：腰圍 64，臀圍 94，大腿圍 60，褲腳圍 58，褲長 95
M：腰圍 68，臀圍 98，大腿圍 62，褲腳圍 60，褲長 96
L：腰圍 72，臀圍 102，大腿圍 64，褲腳圍 62，褲長 97
XL：腰圍 76，臀圍 106，大腿圍 66，褲腳圍 64，褲長 98
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A2").Value = $aText
$ws.Range("B2").Value = $bText
$ws.Rows.Item(2).RowHeight = 205.15

# Row 3
$aText = @'
【YG9631】實拍秋季休閒慵懶褶皺假兩件針織衫(6%羊毛)240925
'@
$bText = @'
肩寬 38，袖長 63，胸圍 76，下擺圍 74，衣長 58
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A3").Value = $aText
$ws.Range("B3").Value = $bText
$ws.Rows.Item(3).RowHeight = 77.650000000000006

# Row 4
$aText = @'
【Y9789】實拍法式復古高腰百褶A字裙240925
'@
$bText = @'
M：腰圍 60，臀圍 110，下擺圍 226，袖長 92
L：腰圍 64，臀圍 114，下擺圍 230，袖長 93
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準
'@
$ws.Range("A4").Value = $aText
$ws.Range("B4").Value = $bText
$ws.Rows.Item(4).RowHeight = 128.65

# Row 5
$aText = @'
【Y82136】實拍韓國寬鬆不規則側開叉斜扣純棉襯衫240925
'@
$bText = @'
M：肩寬 48，袖長 59.5，胸圍 122，下擺圍 130，衣長 71/75
L：肩寬 49，袖長 60.5，胸圍 126，下擺圍 134，衣長 72/76
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準
'@
$ws.Range("A5").Value = $aText
$ws.Range("B5").Value = $bText
$ws.Rows.Item(5).RowHeight = 115.9

# Row 6
$aText = @'
【Y9790】實拍韓國設計款休閒質感衞褲240925
'@
$bText = @'
M：腰圍 66，臀圍 124，大腿圍 76，褲腳圍 58，褲長 102
L：腰圍 70，臀圍 128，大腿圍 78，褲腳圍 60，褲長 103
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準
'@
$ws.Range("A6").Value = $aText
$ws.Range("B6").Value = $bText
$ws.Rows.Item(6).RowHeight = 115.9

# Row 7
$aText = @'
【YA12】實拍軟軟舒服的針織外套+針織背心 套裝240925
'@
$bText = @'
F（背心）：胸圍 72，下擺圍 78，衣長 48
F（外套）：肩寬 60，袖長 46，胸圍 110，下擺圍 110，衣長 76
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。

'@
$ws.Range("A7").Value = $aText
$ws.Range("B7").Value = $bText
$ws.Range("B7").Font.Name = "Arial"
$ws.Range("B7").Font.Size = 10
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("B7").WrapText = $false
$ws.Rows.Item(7).RowHeight = 63.75

# Row 8
$aText = @'
【YP661】實拍慵懶垂感高腰針織闊腿褲240925
'@
$bText = @'
腰圍 61，臀圍 120，大腿圍 88，褲腳圍 134，褲長 94
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。

'@
$ws.Range("A8").Value = $aText
$ws.Range("B8").Value = $bText
$ws.Rows.Item(8).RowHeight = 103.15

# Row 9
$aText = @'
【YG9662】實拍休閒撞色假兩件T恤針織衫240925
'@
$bText = @'
肩寬 36，袖長 61，胸圍 76，下擺圍 76，衣長 55
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A9").Value = $aText
$ws.Range("B9").Value = $bText
$ws.Rows.Item(9).RowHeight = 64.900000000000006

# Row 10
$aText = @'
【YG5162】實拍氣質百搭方領坑條長袖針織衫240925
'@
$bText = @'
肩寬 35，袖長 58，胸圍 66，下擺圍 66，衣長 55
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A10").Value = $aText
$ws.Range("B10").Value = $bText
$ws.Range("B10").Font.Name = "Arial"
$ws.Range("B10").Font.Size = 10
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 64.900000000000006

# Row 11
$aText = @'
【YA8】實拍韓系簡約V領寬鬆針織外套240925
'@
$bText = @'
袖長 67，胸圍 100，下擺圍 86，衣長 52
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A11").Value = $aText
$ws.Range("B11").Value = $bText
$ws.Rows.Item(11).RowHeight = 77.650000000000006

# Row 12
$aText = @'
【Y68106】實拍S~XL高腰顯瘦直筒牛仔褲240925
'@
$bText = @'
S：腰圍 66，臀圍 96，大腿圍 60，褲腳圍 50，褲長 103
M：腰圍 70，臀圍 100，大腿圍 62，褲腳圍 52，褲長 104
L：腰圍 74，臀圍 104，大腿圍 64，褲腳圍 54，褲長 105
XL：腰圍 78，臀圍 108，大腿圍 66，褲腳圍 56，褲長 106
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A12").Value = $aText
$ws.Range("B12").Value = $bText
$ws.Rows.Item(12).RowHeight = 192.4

# Row 13
$aText = @'
【YH9607】實拍S~XL復古高腰直筒牛仔褲240925
'@
$bText = @'
S：腰圍 64，臀圍 92，大腿圍 60，褲腳圍 48，褲長 91
M：腰圍 68，臀圍 96，大腿圍 62，褲腳圍 50，褲長 92
L：腰圍 72，臀圍 100，大腿圍 64，褲腳圍 52，褲長 93
XL：腰圍 76，臀圍 104，大腿圍 66，褲腳圍 54，褲長 94
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A13").Value = $aText
$ws.Range("B13").Value = $bText
$ws.Rows.Item(13).RowHeight = 192.4

# Row 14
$aText = @'
【YP332】實拍韓版半高領軟糯針織T恤洋裝240925
'@
$bText = @'
肩寬 39，袖長 22，胸圍 94，下擺圍 104，衣長 124
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A14").Value = $aText
$ws.Range("B14").Value = $bText
$ws.Rows.Item(14).RowHeight = 77.650000000000006

# Row 15
$aText = @'
【YG9833】實拍半拉鍊假兩件針織T恤(6%羊毛)240925
'@
$bText = @'
肩寬 38，袖長 63，胸圍 70，下擺圍 68，衣長 50/55
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A15").Value = $aText
$ws.Range("B15").Value = $bText
$ws.Rows.Item(15).RowHeight = 77.650000000000006

# Row 16
$aText = @'
【YG736】實拍簡約設計感修身長袖針織T恤(6%羊毛)240925
'@
$bText = @'
肩寬 33，袖長 62，胸圍 72，下擺圍 74，衣長 56
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A16").Value = $aText
$ws.Range("B16").Value = $bText
$ws.Range("A16").Font.Name = "Arial"
$ws.Range("A16").Font.Size = 10
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4108
$ws.Range("A16").WrapText = $true
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 64.900000000000006

# Row 17
$aText = @'
【YA13】實拍氣質針織外套＋背心兩件套240925
'@
$bText = @'
F（背心）：胸圍 68，下擺圍 66，衣長 49
F（外套）：肩寬 57，袖長 54，胸圍 116，下擺圍 114，衣長 58
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A17").Value = $aText
$ws.Range("B17").Value = $bText
$ws.Range("A17").Font.Name = "Arial"
$ws.Range("A17").Font.Size = 10
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4108
$ws.Range("A17").WrapText = $true
$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10
$ws.Range("B17").VerticalAlignment = -4108
$ws.Range("B17").WrapText = $false
$ws.Rows.Item(17).RowHeight = 51

# Row 18
$aText = @'
【YXXT325】實拍簡約粗針針織背心(羊毛25)240925
'@
$bText = @'
肩寬 47，胸圍 100，下擺圍 90，衣長 54
手工平鋪測量，誤差允許在1~3cm左右，具體以實物為準。
'@
$ws.Range("A18").Value = $aText
$ws.Range("B18").Value = $bText
$ws.Range("A18").Font.Name = "Arial"
$ws.Range("A18").Font.Size = 10
$ws.Range("A18").HorizontalAlignment = -4108
$ws.Range("A18").VerticalAlignment = -4108
$ws.Range("A18").WrapText = $true
$ws.Range("B18").Font.Name = "Arial"
$ws.Range("B18").Font.Size = 10
$ws.Range("B18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 64.900000000000006

$null = $ws.Range("A2:B18").Select()
